# Minor fixes to UI
# - Insert a new row at the top of the "Events" data table (row 2) for a
#   "Not Related" event with blank date/detail fields.
# - Renumber the Id column sequentially.
# - Make "Events" the active sheet/tab instead of "PhoneType".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# Insert a new row before the current row 2, shifting the existing rows down.
$ws.Rows(2).Insert(-4121)   # -4121 = xlShiftDown

# Populate the newly inserted row with the "Not Related" placeholder event.
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Not Related"

# Leave the remaining detail columns blank (no value, default/no style).
$ws.Range("C2:E2").ClearContents()
$ws.Range("C2:E2").Style = "Normal"
$ws.Cells.Item(2, 7).ClearContents()
$ws.Cells.Item(2, 7).Style = "Normal"
$ws.Cells.Item(2, 8).ClearContents()
$ws.Cells.Item(2, 8).Style = "Normal"
$ws.Cells.Item(2, 10).ClearContents()
$ws.Cells.Item(2, 10).Style = "Normal"
$ws.Cells.Item(2, 12).ClearContents()
$ws.Cells.Item(2, 12).Style = "Normal"

# Renumber the Id column (A) for the rows that follow.
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(6, 1).Value = 5

# Make the Events sheet the active / selected tab (was previously PhoneType).
$ws.Activate()
$ws.Range("A7").Select()
